$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Editor_Stats")
Write-Host $ws.Name
